$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1. Remove the "TransactionVoucher" sheet entirely (its comments part
#    is dropped automatically along with the sheet by the engine).
# ------------------------------------------------------------------
$wsVoucher = $wb.Worksheets.Item("TransactionVoucher")
$wsVoucher.Delete()

# ------------------------------------------------------------------
# 2. SingleAccount sheet: insert a new "Cheque no" column between
#    "Instrument type" (D) and "Transaction code" (old E, new F),
#    update the existing rows and append a new data row.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SingleAccount")

# insert blank column before E - shifts old E,F to F,G
$ws1.Columns.Item(5).Insert()
$ws1.Columns.Item(5).ColumnWidth = 17.6

# the header comment that used to live on F1 ("Transaction code")
# now belongs on G1 - recreate it there with the same text
$oldCommentText = $ws1.Range("F1").Comment.Text()
$ws1.Range("F1").Comment.Delete()
$ws1.Range("G1").AddComment($oldCommentText)

# header row
$ws1.Range("E1").Value = "Cheque no"

# row 2 : 111423-Trial Entry / SB 55-Trial Entry, Credit -> Debit,
# new Cheque no value
$ws1.Range("C2").Value = "Debit"
$ws1.Range("E2").Value = 605488

# row 3 : replace the old duplicate "Trial Entry" row with the
# "111974-Pranav Parab" / Withdrawal Slip row (format/style is already
# correct because the column insert shifted the existing cells)
$ws1.Range("A3").Value = "111974-Pranav Parab"
$ws1.Range("B3").Value = "SB 116-Pranav Parab"
$ws1.Range("C3").Value = "Debit"
$ws1.Range("D3").Value = "Withdrawal Slip"
$ws1.Range("E3").Value = ""
$ws1.Range("F3").Value = "SB"
$ws1.Range("G3").Value = 8000

# row 4 : brand new "111851-abcd ijkl" / Reciept row - clone row 3's
# formatting (style + row height) then overwrite with the new values
$ws1.Range("A3:G3").Copy()
$ws1.Range("A4:G4").PasteSpecial(-4122)
$ws1.Rows.Item(4).RowHeight = 53.25

$ws1.Range("A4").Value = "111851-abcd ijkl"
$ws1.Range("B4").Value = "SB 101-abcd ijkl"
$ws1.Range("C4").Value = "Credit"
$ws1.Range("D4").Value = "Reciept"
$ws1.Range("E4").Value = ""
$ws1.Range("F4").Value = "CA"
$ws1.Range("G4").Value = 10000

$ws1.Activate()
$ws1.Range("E2").Select()

# ------------------------------------------------------------------
# 3. MultipleAccount sheet: selection only changes (K1 -> E2); the
#    underlying cell values are untouched.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MultipleAccount")
$ws2.Activate()
$ws2.Range("E2").Select()
